# Commit message: "rename Collection to CRF in tabs"
#
# Rename the sheet from Collection_DS to CRF_DS (the worksheet tab name).
# Excel automatically keeps the _xlnm._FilterDatabase defined name (and any
# formulas referencing the sheet) in sync with the new sheet name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab: Collection_DS -> CRF_DS
$ws.Name = "CRF_DS"

# Update the active selection on the frozen (bottom-left) pane from R1 to K5,
# matching the saved cursor position recorded in the workbook.
$ws.Range("K5").Select()
